$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Forecasted Consumption (MW)" values for A2:A97 (96 rows)
$consumption = @(
    5270,5220,5180,5140,5110,5080,5070,5060,5050,5050,5050,5050,
    5050,5060,5070,5080,5090,5120,5160,5220,5280,5360,5420,5510,
    5640,5770,5860,5910,5920,5920,5910,5880,5800,5700,5610,5520,
    5450,5370,5310,5250,5150,5090,5060,5040,5020,5020,5020,5020,
    5030,5040,5050,5080,5100,5130,5160,5220,5280,5340,5390,5460,
    5530,5620,5720,5830,5910,6020,6110,6210,6300,6400,6500,6580,
    6660,6730,6780,6850,6910,6960,6990,7000,7000,6990,6970,6900,
    6720,6550,6390,6230,6090,5950,5840,5720,5650,5590,5540,5470
)

for ($i = 0; $i -lt $consumption.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $consumption[$i]
    # Shift the Timestamp column forward by 14 days (45875 -> 45889 base)
    $oldTimestamp = $ws.Cells.Item($row, 2).Value2
    $ws.Cells.Item($row, 2).Value = $oldTimestamp + 14
}
